$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: seed new Batter/Bowler names in the exact order they were
# first introduced in the source workbook, so new shared-string entries
# land in the same order as the authored file.
$ws.Cells.Item(82, 14).Value2 = 'RM Patidar'
$ws.Cells.Item(82, 15).Value2 = 'JJ Bumrah'
$ws.Cells.Item(84, 14).Value2 = 'MK Lomror'
$ws.Cells.Item(86, 14).Value2 = 'SA Yadav'
$ws.Cells.Item(86, 15).Value2 = 'V Vyshak'
$ws.Cells.Item(83, 14).Value2 = 'F du Plessis'
$ws.Cells.Item(83, 15).Value2 = 'A Madhwal'

# Step 2: fill in the remaining cells for the new match rows (82-86),
# covering every column; N/O for rows 82/83/84/86 were already seeded above.

# Row 82
$ws.Cells.Item(82, 1).Value2 = 25
$ws.Cells.Item(82, 2).Value2 = 'MI'
$ws.Cells.Item(82, 3).Value2 = 'RCB'
$ws.Cells.Item(82, 4).Value2 = 1
$ws.Cells.Item(82, 5).Value2 = 'RCB'
$ws.Cells.Item(82, 6).Value2 = 'MI'
$ws.Cells.Item(82, 7).Value2 = 11
$ws.Cells.Item(82, 8).Value2 = 'MI'
$ws.Cells.Item(82, 9).Value2 = 'Nitin Menon'
$ws.Cells.Item(82, 10).Value2 = 'NM'
$ws.Cells.Item(82, 11).Value2 = 'Wicket'
$ws.Cells.Item(82, 12).Value2 = 'Not Out'
$ws.Cells.Item(82, 13).Value2 = 'Not Out'
$ws.Cells.Item(82, 16).Value2 = 'Unsuccessful'
$ws.Cells.Item(82, 17).Value2 = 'No'

# Row 83
$ws.Cells.Item(83, 1).Value2 = 25
$ws.Cells.Item(83, 2).Value2 = 'MI'
$ws.Cells.Item(83, 3).Value2 = 'RCB'
$ws.Cells.Item(83, 4).Value2 = 1
$ws.Cells.Item(83, 5).Value2 = 'RCB'
$ws.Cells.Item(83, 6).Value2 = 'MI'
$ws.Cells.Item(83, 7).Value2 = 16
$ws.Cells.Item(83, 8).Value2 = 'MI'
$ws.Cells.Item(83, 9).Value2 = 'VA Kulkarni'
$ws.Cells.Item(83, 10).Value2 = 'VAK'
$ws.Cells.Item(83, 11).Value2 = 'Wide'
$ws.Cells.Item(83, 12).Value2 = 'Called'
$ws.Cells.Item(83, 13).Value2 = 'Called'
$ws.Cells.Item(83, 16).Value2 = 'Unsuccessful'
$ws.Cells.Item(83, 17).Value2 = 'No'

# Row 84
$ws.Cells.Item(84, 1).Value2 = 25
$ws.Cells.Item(84, 2).Value2 = 'MI'
$ws.Cells.Item(84, 3).Value2 = 'RCB'
$ws.Cells.Item(84, 4).Value2 = 1
$ws.Cells.Item(84, 5).Value2 = 'RCB'
$ws.Cells.Item(84, 6).Value2 = 'MI'
$ws.Cells.Item(84, 7).Value2 = 17
$ws.Cells.Item(84, 8).Value2 = 'RCB'
$ws.Cells.Item(84, 9).Value2 = 'Nitin Menon'
$ws.Cells.Item(84, 10).Value2 = 'NM'
$ws.Cells.Item(84, 11).Value2 = 'Wicket'
$ws.Cells.Item(84, 12).Value2 = 'Out'
$ws.Cells.Item(84, 13).Value2 = 'Out'
$ws.Cells.Item(84, 15).Value2 = 'JJ Bumrah'
$ws.Cells.Item(84, 16).Value2 = 'Unsuccessful'
$ws.Cells.Item(84, 17).Value2 = 'Yes'

# Row 85
$ws.Cells.Item(85, 1).Value2 = 25
$ws.Cells.Item(85, 2).Value2 = 'MI'
$ws.Cells.Item(85, 3).Value2 = 'RCB'
$ws.Cells.Item(85, 4).Value2 = 1
$ws.Cells.Item(85, 5).Value2 = 'RCB'
$ws.Cells.Item(85, 6).Value2 = 'MI'
$ws.Cells.Item(85, 7).Value2 = 20
$ws.Cells.Item(85, 8).Value2 = 'RCB'
$ws.Cells.Item(85, 9).Value2 = 'VA Kulkarni'
$ws.Cells.Item(85, 10).Value2 = 'VAK'
$ws.Cells.Item(85, 11).Value2 = 'NoBall'
$ws.Cells.Item(85, 12).Value2 = 'Not Called'
$ws.Cells.Item(85, 13).Value2 = 'Not Called'
$ws.Cells.Item(85, 14).Value2 = 'KD Karthik'
$ws.Cells.Item(85, 15).Value2 = 'A Madhwal'
$ws.Cells.Item(85, 16).Value2 = 'Unsuccessful'
$ws.Cells.Item(85, 17).Value2 = 'No'

# Row 86
$ws.Cells.Item(86, 1).Value2 = 25
$ws.Cells.Item(86, 2).Value2 = 'MI'
$ws.Cells.Item(86, 3).Value2 = 'RCB'
$ws.Cells.Item(86, 4).Value2 = 2
$ws.Cells.Item(86, 5).Value2 = 'MI'
$ws.Cells.Item(86, 6).Value2 = 'RCB'
$ws.Cells.Item(86, 7).Value2 = 14
$ws.Cells.Item(86, 8).Value2 = 'MI'
$ws.Cells.Item(86, 9).Value2 = 'Nitin Menon'
$ws.Cells.Item(86, 10).Value2 = 'NM'
$ws.Cells.Item(86, 11).Value2 = 'NoBall'
$ws.Cells.Item(86, 12).Value2 = 'Not Called'
$ws.Cells.Item(86, 13).Value2 = 'Not Called'
$ws.Cells.Item(86, 16).Value2 = 'Unsuccessful'
$ws.Cells.Item(86, 17).Value2 = 'No'

# View state: reflect the new scroll position / zoom / selection used while
# reviewing the freshly-added match rows.
$win = $excel.ActiveWindow
$win.Zoom = 110
$ws.Range("N85").Select()
